$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 data
$ws.Range("F11").Value = 45626
$ws.Range("F11").NumberFormat = "mm-dd-yy"
$ws.Range("G11").Value = 0.47916666666666669
$ws.Range("I11").Formula = "=H11-G11"

# Update the selection to H11 (matches diff's <selection activeCell="H11" sqref="H11"/>)
$ws.Range("H11").Select()
